$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CALC + CDP re-correction: update Index/Name/sem/CGPA cells for the affected rows
# (rows 42-53 and 85-89) to match the corrected ranking/recalculation.

# Row 42
$ws.Cells.Item(42, 2).Value = 230045
$ws.Cells.Item(42, 3).Value = 'ANTHONY C.S.B.'
$ws.Cells.Item(42, 6).Value = 3.704
$ws.Cells.Item(42, 7).Value = 3.886

# Row 43
$ws.Cells.Item(43, 2).Value = 230038
$ws.Cells.Item(43, 3).Value = 'AMARATHUNGE A.M.N.L.'
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(43, 6).Value = 3.691

# Row 44
$ws.Cells.Item(44, 2).Value = 230130
$ws.Cells.Item(44, 3).Value = 'DESHAN W.U.'
$ws.Cells.Item(44, 4).Value = 4
$ws.Cells.Item(44, 5).Value = 3.96
$ws.Cells.Item(44, 6).Value = 3.714
$ws.Cells.Item(44, 7).Value = 3.881

# Row 45
$ws.Cells.Item(45, 2).Value = 230212
$ws.Cells.Item(45, 3).Value = 'GUNASEKARA L.U.A.'
$ws.Cells.Item(45, 4).Value = 3.957
$ws.Cells.Item(45, 5).Value = 3.882
$ws.Cells.Item(45, 6).Value = 3.821

# Row 46
$ws.Cells.Item(46, 2).Value = 230321
$ws.Cells.Item(46, 3).Value = 'KARUNANAYAKE A.H.D.'
$ws.Cells.Item(46, 4).Value = 4
$ws.Cells.Item(46, 5).Value = 3.947
$ws.Cells.Item(46, 6).Value = 3.73
$ws.Cells.Item(46, 7).Value = 3.876

# Row 47
$ws.Cells.Item(47, 2).Value = 230300
$ws.Cells.Item(47, 3).Value = 'JAYAWEERA N.S.'
$ws.Cells.Item(47, 6).Value = 3.713
$ws.Cells.Item(47, 7).Value = 3.875

# Row 48
$ws.Cells.Item(48, 2).Value = 230145
$ws.Cells.Item(48, 3).Value = 'DILHAN W.A.'
$ws.Cells.Item(48, 5).Value = 4
$ws.Cells.Item(48, 6).Value = 3.704
$ws.Cells.Item(48, 7).Value = 3.871

# Row 49
$ws.Cells.Item(49, 2).Value = 230477
$ws.Cells.Item(49, 3).Value = 'PERERA H.A.J.I.'
$ws.Cells.Item(49, 4).Value = 3.935
$ws.Cells.Item(49, 5).Value = 3.817
$ws.Cells.Item(49, 6).Value = 3.873
$ws.Cells.Item(49, 7).Value = 3.866

# Row 50
$ws.Cells.Item(50, 2).Value = 230613
$ws.Cells.Item(50, 3).Value = 'SHEHAN M.N.N.'
$ws.Cells.Item(50, 5).Value = 3.947
$ws.Cells.Item(50, 6).Value = 3.699
$ws.Cells.Item(50, 7).Value = 3.865

# Row 51
$ws.Cells.Item(51, 2).Value = 230058
$ws.Cells.Item(51, 3).Value = 'AROSHANA H.A.P.'
$ws.Cells.Item(51, 4).Value = 4
$ws.Cells.Item(51, 5).Value = 3.908
$ws.Cells.Item(51, 6).Value = 3.713

# Row 52
$ws.Cells.Item(52, 2).Value = 230697
$ws.Cells.Item(52, 3).Value = 'WEERASINGHE J.A.H.R.'
$ws.Cells.Item(52, 4).Value = 3.957
$ws.Cells.Item(52, 5).Value = 3.96
$ws.Cells.Item(52, 6).Value = 3.671
$ws.Cells.Item(52, 7).Value = 3.855

# Row 53
$ws.Cells.Item(53, 2).Value = 230211
$ws.Cells.Item(53, 3).Value = 'GUNASEKARA K.S.'
$ws.Cells.Item(53, 5).Value = 3.895
$ws.Cells.Item(53, 6).Value = 3.721
$ws.Cells.Item(53, 7).Value = 3.853

# Row 85
$ws.Cells.Item(85, 2).Value = 230495
$ws.Cells.Item(85, 3).Value = 'PRABHARSHA H.W.D.'
$ws.Cells.Item(85, 4).Value = 3.85
$ws.Cells.Item(85, 5).Value = 3.869
$ws.Cells.Item(85, 6).Value = 3.443
$ws.Cells.Item(85, 7).Value = 3.701

# Row 86
$ws.Cells.Item(86, 2).Value = 230444
$ws.Cells.Item(86, 3).Value = 'NIRMANI W.T.'
$ws.Cells.Item(86, 4).Value = 3.785
$ws.Cells.Item(86, 5).Value = 3.596
$ws.Cells.Item(86, 6).Value = 3.757
$ws.Cells.Item(86, 7).Value = 3.693

# Row 87
$ws.Cells.Item(87, 2).Value = 230261
$ws.Cells.Item(87, 3).Value = 'INDUWARA M.L.A.S.'
$ws.Cells.Item(87, 4).Value = 4
$ws.Cells.Item(87, 5).Value = 3.747
$ws.Cells.Item(87, 6).Value = 3.458

# Row 88
$ws.Cells.Item(88, 2).Value = 230375
$ws.Cells.Item(88, 3).Value = 'LENMINI B.L.W.'
$ws.Cells.Item(88, 4).Value = 3.85
$ws.Cells.Item(88, 5).Value = 3.686
$ws.Cells.Item(88, 6).Value = 3.599
$ws.Cells.Item(88, 7).Value = 3.691

# Row 89
$ws.Cells.Item(89, 2).Value = 230527
$ws.Cells.Item(89, 3).Value = 'RANAWAKA R.A.G.K.'
$ws.Cells.Item(89, 5).Value = 3.83
$ws.Cells.Item(89, 6).Value = 3.339
$ws.Cells.Item(89, 7).Value = 3.681
